# 391-RBI-EI-DB-DL-REC-NON-RNI-CTRFD-SAR-MD-TR-1-EarlyRePayment-Loanproduct4.xlsx
# code refactoring and loan accounting and charges added

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- shortname is now numeric (391) instead of the text "kar5" ---
$ws.Range("B3").Value = 391

# --- nominalinterestratedefault changes from 12 to 1 ---
$ws.Range("B11").Value = 1

# --- new loan accounting / charges rows (29-40), copying the format of an
#     existing data row (row 10) so the cells reuse the workbook's existing
#     styles instead of minting new ones ---
$ws.Range("A10:B10").Copy()
$ws.Range("A29:B40").PasteSpecial(-4122)

# column B (values) was filled in first in the original edit, so write it
# first here too -- this keeps the regenerated sharedStrings.xml in the same
# order as the authored workbook.
$ws.Range("B29").Value = "Cash"
$ws.Range("B30").Value = "Loan portfolio "
$ws.Range("B31").Value = "Interest Receivable "
$ws.Range("B32").Value = "Penalties Receivable "
$ws.Range("B33").Value = "Transfer in Suspence "
$ws.Range("B34").Value = "Fees Receivable"
$ws.Range("B35").Value = "Income from interest"
$ws.Range("B36").Value = "Income from penalties"
$ws.Range("B37").Value = "Income from fees"
$ws.Range("B38").Value = "Income from recovery repayments"
$ws.Range("B39").Value = "Losses Writtenoff "
$ws.Range("B40").Value = "Overpayment Liability"

# then column A (labels)
$ws.Range("A29").Value = "fundsource"
$ws.Range("A30").Value = "loanprotfolio"
$ws.Range("A31").Value = "interestreceivable"
$ws.Range("A32").Value = "penaltiesreceivable"
$ws.Range("A33").Value = "transferinsuspense"
$ws.Range("A34").Value = "feesreceivable"
$ws.Range("A35").Value = "incomefrominterest"
$ws.Range("A36").Value = "incomefrompenalties"
$ws.Range("A37").Value = "incomefromfees"
$ws.Range("A38").Value = "incomefromrecoveryrepayments"
$ws.Range("A39").Value = "loseswrittenoff"
$ws.Range("A40").Value = "overpaymentliability"

# --- column B is now wider to fit the longer accounting descriptions ---
$ws.Columns.Item(2).ColumnWidth = 61.7109375

# --- move the active selection like the source workbook's last save ---
$ws.Range("A68").Select()

# --- sheet2 mirrors the product name in B1 (same shared string as sheet1!B1) ---
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("B1").Value = "391-RBI-EI-DB-DL-REC-NON-RNI-CTRFD-SAR-MD-TR-1-EarlyRePayment"
